$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels in row 1 to include units (and fix accent on "Posición")
$ws.Range("A1").Value = "Posición [cm]"
$ws.Range("B1").Value = "Tiempo [s]"

# Update the selected/active cell to C3 (as saved in the workbook)
$ws.Range("C3").Select()
